# Updated cryptocurrency price/volume data to reflect the latest scrape.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Rows 10-18 and 42-43 were
# reordered because the ranking of several tokens changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    # D column holds numeric-looking values that must stay stored as text
    # (matching the original inlineStr / shared-string text cells), so we
    # force the cell to Text format while assigning, then clear the
    # formatting again so no extra number-format style lingers on the cell.
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

function Set-PlainValue($cell, $val) {
    $ws.Range($cell).Value = $val
}

Set-TextValue 'D2' '243.01'
Set-TextValue 'D3' '23.35'
Set-TextValue 'D4' '5.644'
Set-TextValue 'D5' '0.05810'
Set-TextValue 'D6' '3.412'
Set-TextValue 'D7' '6.467'
Set-TextValue 'D8' '1.320'
Set-TextValue 'D9' '0.7980'
Set-PlainValue 'B10' 'WazirX'
Set-PlainValue 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D10' '0.1458'
Set-PlainValue 'E10' '9WazirXWRX'
Set-PlainValue 'B11' 'MandalaExchangeToken'
Set-PlainValue 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D11' '0.07630'
Set-PlainValue 'E11' '10MandalaExchangeTokenMDX'
Set-PlainValue 'B12' 'LiechtensteinCryptoassetsExchange'
Set-PlainValue 'C12' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D12' '0.03230'
Set-PlainValue 'E12' '11LiechtensteinCryptoassetsExchangeLCX'
Set-PlainValue 'B13' 'BitrueCoin'
Set-PlainValue 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.02963'
Set-PlainValue 'E13' '12BitrueCoinBTR'
Set-PlainValue 'B14' 'BitMartToken'
Set-PlainValue 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.09241'
Set-PlainValue 'E14' '13BitMartTokenBMX'
Set-PlainValue 'B15' 'BitForexToken'
Set-PlainValue 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D15' '0.001663'
Set-PlainValue 'E15' '14BitForexTokenBF'
Set-PlainValue 'B16' 'MCDex'
Set-PlainValue 'C16' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D16' '3.327'
Set-PlainValue 'E16' '15MCDexMCB'
Set-PlainValue 'B17' 'CoinExToken'
Set-PlainValue 'C17' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue 'D17' '0.04749'
Set-PlainValue 'E17' '16CoinExTokenCET'
Set-PlainValue 'B18' 'One'
Set-PlainValue 'C18' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 'D18' '0.0005992'
Set-PlainValue 'E18' '17OneONE'
Set-TextValue 'D20' '0.005461'
Set-TextValue 'D21' '0.001067'
Set-TextValue 'D22' '0.0001500'
Set-TextValue 'D24' '2.192'
Set-TextValue 'D26' '0.1241'
Set-TextValue 'D27' '0.001000'
Set-TextValue 'D40' '0.04279'
Set-TextValue 'D41' '0.007140'
Set-PlainValue 'B42' 'BKEXToken'
Set-PlainValue 'C42' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D42' '0.1055'
Set-PlainValue 'E42' '41BKEXTokenBKK'
Set-PlainValue 'B43' 'CEJI'
Set-PlainValue 'C43' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue 'D43' '0.003170'
Set-PlainValue 'E43' '42CEJICEJI'
Set-TextValue 'D44' '0.009539'
Set-TextValue 'D46' '0.00005440'
Set-TextValue 'D47' '0.00000000750'
Set-TextValue 'D48' '0.7855'
Set-TextValue 'D49' '0.1024'
Set-PlainValue 'E49' '48BOLOBOLOBestin24h'
Set-TextValue 'D50' '0.00002101'
Set-TextValue 'D51' '0.01010'
